# Fruta / hortaliza, semanal
# Insert 2 new weekly rows (new date 2022-03-17 / serial 44637) above the
# existing row 22 ("Especial" and "Primera" quality entries), pushing all
# subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 22 - existing rows 22:27 shift down to 24:29.
$ws.Rows("22:23").Insert()

# New row 22: Tuna "Especial"
$ws.Range("A22").Value = 8
$ws.Range("B22").Value = "Terminal La Palmera de La Serena"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 44637
$ws.Range("D22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100107
$ws.Range("H22").Value = "Otros"
$ws.Range("I22").Value = 100107011
$ws.Range("J22").Value = "Tuna"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Especial"
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 14000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 14500
$ws.Range("Q22").Value = "$/caja 18 kilos"
$ws.Range("R22").Value = "Provincia de Limarí"
$ws.Range("S22").Value = 806
$ws.Range("T22").Value = 18

# New row 23: Tuna "Primera"
$ws.Range("A23").Value = 8
$ws.Range("B23").Value = "Terminal La Palmera de La Serena"
$ws.Range("C23").Value = "Coquimbo"
$ws.Range("D23").Value = 44637
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100107
$ws.Range("H23").Value = "Otros"
$ws.Range("I23").Value = 100107011
$ws.Range("J23").Value = "Tuna"
$ws.Range("K23").Value = "Sin especificar"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 240
$ws.Range("N23").Value = 10000
$ws.Range("O23").Value = 11000
$ws.Range("P23").Value = 10500
$ws.Range("Q23").Value = "$/caja 18 kilos"
$ws.Range("R23").Value = "Provincia de Limarí"
$ws.Range("S23").Value = 583
$ws.Range("T23").Value = 18
